# "ajustes gerais na documentação"
# 1) Refresh the cached text of the "datetimeFigureOut" date placeholder
#    field (slide master + every slide layout) from 24/05/2020 to 15/11/2020.
# 2) Remove the now-obsolete "v1"/"v2" architecture-diagram callouts and the
#    duplicated PostgreSQL flowchart/connector shapes from slides 1 and 2.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = "15/11/2020"
        }
    }
}

# --- Slide master ---
$masterShapes = $p.SlideMaster.Shapes
Update-DatePlaceholder $masterShapes

# --- Every slide layout ---
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layoutShapes = $layouts.Item($li).Shapes
    Update-DatePlaceholder $layoutShapes
}

function Remove-NamedShape {
    param($slide, $name)

    for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            $sh.Delete()
            return
        }
    }
}

# --- Slide 1: drop the "v2" diagram's extra PostgreSQL disk + callouts ---
$s1 = $p.Slides.Item(1)
Remove-NamedShape $s1 "CaixaDeTexto 42"
Remove-NamedShape $s1 "Fluxograma: Disco Magnético 26"
Remove-NamedShape $s1 "CaixaDeTexto 29"
Remove-NamedShape $s1 "Conector de Seta Reta 30"
Remove-NamedShape $s1 "Conector de Seta Reta 31"

# --- Slide 2: drop the "v1" diagram's extra PostgreSQL disk + callouts ---
$s2 = $p.Slides.Item(2)
Remove-NamedShape $s2 "CaixaDeTexto 42"
Remove-NamedShape $s2 "Agrupar 65"
Remove-NamedShape $s2 "Conector de Seta Reta 60"
Remove-NamedShape $s2 "Conector de Seta Reta 61"
